$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Productdata")
$ws.Range("C2").Value = 1
$ws.Range("E2").Value = 1008.888888888889
$ws.Range("C3").Value = 0
$ws.Range("E3").Value = 205.6666666666666
$ws.Range("C4").Value = 0
$ws.Range("E4").Value = 102.7777777777778
$ws.Range("C5").Value = 0
$ws.Range("E5").Value = 410.4444444444443
$ws.Range("C6").Value = 0
$ws.Range("E6").Value = 446.3583333333332
$ws.Range("C7").Value = 60001
$ws.Range("E7").Value = 205.6666666666666
$ws.Range("C8").Value = 30001
$ws.Range("E8").Value = 102.7777777777778
$ws.Range("C9").Value = 120001
$ws.Range("E9").Value = 820.8888888888887
$ws.Range("C10").Value = 10501
$ws.Range("E10").Value = 35.91388888888888
$ws.Range("C11").Value = 0
$ws.Range("E11").Value = 179
$ws.Range("C12").Value = 0
$ws.Range("E12").Value = 89.44444444444443
$ws.Range("C13").Value = 0
$ws.Range("E13").Value = 357.111111111111
$ws.Range("C14").Value = 0
$ws.Range("E14").Value = 388.3583333333333

$wsCap = $wb.Worksheets.Item("Capacity")
$wsCap.Range("B2").Value = 227000
$wsCap.Range("B3").Value = 40000
$wsCap.Range("B4").Value = 20000
$wsCap.Range("B5").Value = 160000
$wsCap.Range("B6").Value = 261000
$wsCap.Range("B7").Value = 80000
$wsCap.Range("B8").Value = 80000
$wsCap.Range("B9").Value = 160000
$wsCap.Range("B10").Value = 35000
$wsCap.Range("B11").Value = 200000
$wsCap.Range("B12").Value = 80000
$wsCap.Range("B13").Value = 400000
$wsCap.Range("B14").Value = 435000

$wsProc = $wb.Worksheets.Item("ProcessingTime")
$wsProc.Range("B2").Value = 1
$wsProc.Range("C3").Value = 1
$wsProc.Range("D4").Value = 1
$wsProc.Range("E5").Value = 2
$wsProc.Range("F6").Value = 3
$wsProc.Range("G7").Value = 2
$wsProc.Range("H8").Value = 4
$wsProc.Range("J10").Value = 5
$wsProc.Range("K11").Value = 5
$wsProc.Range("M13").Value = 5
$wsProc.Range("N14").Value = 5
